$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new values look like pure numbers,
# so Excel stores them as text (matching the workbook convention) instead of
# auto-converting them to numeric values.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row value updates (Price / Volume(1h) columns)
$ws.Range("D2").Value = "42.118.25"
$ws.Range("E2").Value = "  +0.60%  "
$ws.Range("D3").Value = "2.277.61"
$ws.Range("E3").Value = "  +0.21%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "154.95"
$ws.Range("E5").Value = "  +15,471.77%  "
$ws.Range("D6").Value = "305.03"
$ws.Range("E6").Value = "  +0.30%  "
$ws.Range("D7").Value = "94.02"
$ws.Range("E7").Value = "  +0.81%  "
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "0.490"
$ws.Range("E10").Value = "  +0.50%  "
$ws.Range("D11").Value = "33.97"
$ws.Range("E11").Value = "  +3.86%  "
$ws.Range("D12").Value = "0.0803"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("D14").Value = "6.66"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D16").Value = "14.36"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").Value = "2.270.14"
$ws.Range("E17").Value = "  -1.60%  "
$ws.Range("D18").Value = "0.792"
$ws.Range("E18").Value = "  +3.88%  "
$ws.Range("D19").Value = "42.037.42"
$ws.Range("E19").Value = "  +0.59%  "
$ws.Range("D20").Value = "12.78"
$ws.Range("E20").Value = "  +4.33%  "
$ws.Range("D21").Value = "0.0₃0917"
$ws.Range("E21").Value = "  +0.95%  "
$ws.Range("D22").Value = "5.99"
$ws.Range("E22").Value = "  +0.63%  "
$ws.Range("D23").Value = "68.00"
$ws.Range("E23").Value = "  +0.95%  "
$ws.Range("D24").Value = "243.71"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  +0.77%  "
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").Value = "24.08"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("D29").Value = "35.88"
$ws.Range("E29").Value = "  +5.20%  "
$ws.Range("D30").Value = "9.69"
$ws.Range("E30").Value = "  +0.63%  "
$ws.Range("E31").Value = "  +1.18%  "
$ws.Range("D32").Value = "160.64"
$ws.Range("E32").Value = "  +1.43%  "
$ws.Range("D33").Value = "5.34"
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("D35").Value = "0.0754"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("E36").Value = "  +0.72%  "
$ws.Range("E37").Value = "  +3.90%  "
$ws.Range("D38").Value = "16.98"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("D42").Value = "4.21"
$ws.Range("E42").Value = "  +7.16%  "
$ws.Range("D45").Value = "2.27"
$ws.Range("E45").Value = "  +11.33%  "
$ws.Range("E46").Value = "  +1.51%  "
$ws.Range("D47").Value = "10.23"
$ws.Range("E47").Value = "  -2.17%  "
$ws.Range("D48").Value = "2.92"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").Value = "53.46"
$ws.Range("E49").Value = "  +3.25%  "
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").Value = "72.14"
$ws.Range("E51").Value = "  -1.48%  "

# Row 43/44 swap: EnergySwap now ranks ahead of Maker
$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "19.92"
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.021.75"
$ws.Range("E44").Value = "  -2.46%  "
